$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I0 and IF columns, matching the style used by existing headers (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for the new I (I0) and J (IF) columns, rows 2-17
$data = @{
    2  = @(4, 5)
    3  = @(6, 7)
    4  = @(8, 8)
    5  = @(5, 6)
    6  = @(8, 8)
    7  = @(5, 5)
    8  = @(7, 7)
    9  = @(7, 7)
    10 = @(5, 6)
    11 = @(9, 9)
    12 = @(5, 6)
    13 = @(9, 9)
    14 = @(9, 9)
    15 = @(7, 7)
    16 = @(8, 8)
    17 = @(4, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
